$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# xlPasteFormats - used to clone the existing date-column formatting (border,
# bold font, centered alignment, custom date number format) onto new rows.
$xlPasteFormats = -4122

# New daily rows to append (aggiornamento fino a 1/09/2021): dates
# 2021-08-24 .. 2021-09-01 as serials, followed by the per-comune counts
# in column order A (date) .. AX.
$newRows = @(
    @(44432,7,6,2,61,23,8,8,1,2,8,10,12,21,0,4,3,15,5,5,17,190,0,10,11,12,2,4,2,4,11,7,29,11,10,6,12,20,18,0,25,631,2,7,0,1,2,7,1,9),
    @(44433,8,8,2,54,23,8,8,1,2,6,10,14,20,0,4,3,16,6,5,17,199,0,12,12,14,2,4,2,3,11,7,27,12,12,6,10,19,17,0,22,633,2,7,0,1,2,5,1,9),
    @(44434,7,7,2,61,24,7,8,1,2,3,10,13,23,0,4,3,17,4,6,17,179,0,12,12,10,2,4,2,4,11,6,26,10,11,5,12,17,15,0,25,604,2,5,0,1,2,2,2,8),
    @(44435,5,8,1,67,26,7,10,3,1,2,8,13,20,0,2,2,17,4,5,16,186,0,11,11,6,2,3,3,6,15,10,31,9,11,4,9,20,21,0,26,622,2,7,0,1,3,1,3,4),
    @(44436,4,11,1,72,24,6,8,3,1,3,9,13,20,0,2,2,15,4,6,16,173,0,11,10,10,2,1,4,8,17,10,30,11,15,5,6,11,21,0,29,627,2,9,0,0,3,6,7,6),
    @(44437,11,9,6,70,35,5,7,3,1,4,14,19,18,0,0,1,14,4,4,14,188,0,14,10,18,1,1,5,7,13,10,19,15,13,5,6,11,30,0,28,665,2,9,0,0,2,6,7,6),
    @(44438,12,12,6,107,38,6,9,3,1,2,14,20,16,0,0,0,18,4,4,15,183,2,16,7,22,3,1,11,14,13,12,27,16,16,3,5,11,30,1,29,740,2,8,0,0,2,6,7,6),
    @(44439,13,11,5,113,38,6,9,5,2,5,24,30,15,0,0,0,15,4,2,14,169,2,16,6,25,3,1,11,15,12,13,29,12,15,2,5,10,33,1,35,755,2,7,0,0,2,6,7,5),
    @(44440,13,9,5,114,37,5,9,6,2,5,24,28,15,0,0,0,13,3,2,14,161,2,15,5,23,3,1,11,15,12,13,31,12,13,2,5,10,34,1,34,736,2,7,0,0,2,6,7,5)
)

$startRow = 358
$lastDataRow = 357

for ($r = 0; $r -lt $newRows.Length; $r++) {
    $targetRow = $startRow + $r
    $rowdata = $newRows[$r]

    # Copy the date-column style (border/bold/center + date number format)
    # from the last existing data row onto the new date cell.
    $ws.Range("A" + $lastDataRow).Copy()
    $ws.Range("A" + $targetRow).PasteSpecial($xlPasteFormats)

    for ($c = 0; $c -lt $rowdata.Length; $c++) {
        $ws.Cells.Item($targetRow, $c + 1).Value = $rowdata[$c]
    }
}